$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2: update invoice number / quantity / emission date
$ws.Range("A2").Value = "F002"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 45767

# Row 6: product name and unit price change (formula in D6 recalculates automatically)
$ws.Range("A6").Value = "Cosa1"
$ws.Range("C6").Value = 20000

# Update the selected cell to reflect the saved workbook view
$ws.Range("C7").Select()
